# Trade #57 closed at 2026-02-17 08:47:50 - unknown UNKNOWN +0.000%
#
# Adds a new closed MarketMaking trade row (trade #57 / row 58) to the
# "All Trades" and "MarketMaking" sheets, and updates the aggregate
# statistics on the "Summary" and "Strategy Status" sheets accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "Summary" sheet aggregate metrics
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.6     # Current Capital
$summary.Range("B4").Value = -0.4       # Total P&L $
$summary.Range("B5").Value = -0.14      # Total P&L %
$summary.Range("B6").Value = 57         # Total Trades
$summary.Range("B8").Value = 25         # Losing Trades
$summary.Range("B9").Value = 38.6       # Win Rate %

# ---------------------------------------------------------------------
# 2) Update the "Strategy Status" sheet for the MarketMaking strategy
#    (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.59999999999999   # Capital
$status.Range("D4").Value = 57                  # Trades
$status.Range("E4").Value = -0.4                # P&L $
$status.Range("F4").Value = -0.4                # P&L %
$status.Range("G4").Value = 38.6                # Win Rate %

# ---------------------------------------------------------------------
# 3) Append the new trade row (row 58) to both "All Trades" and
#    "MarketMaking" sheets
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A58").Value = 57

    # Force the date/time-looking strings to be stored as plain text,
    # matching the original sheet's inline-string cells, instead of
    # letting Excel auto-convert them to date/time serial numbers.
    $ws.Range("B58").NumberFormat = "@"
    $ws.Range("B58").Value = "2026-02-17"

    $ws.Range("C58").NumberFormat = "@"
    $ws.Range("C58").Value = "08:47:44"

    $ws.Range("D58").Value = "MarketMaking"
    $ws.Range("E58").Value = "DOWN"
    $ws.Range("F58").Value = 0.34
    $ws.Range("G58").Value = 0.25
    $ws.Range("H58").Value = "CLOSED"
    $ws.Range("I58").Value = -26.4706
    $ws.Range("J58").Value = -0.09
    $ws.Range("K58").Value = 99.59999999999999
    $ws.Range("L58").Value = 0
    $ws.Range("M58").Value = 0
    $ws.Range("N58").Value = 0.6
    $ws.Range("O58").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P58").Value = "early_exit"
    $ws.Range("Q58").Value = 0.13
}
